$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 14: recompute K14
$ws.Range("K14").Value = -0.2254024683979639

# Row 15: recompute J15, add K15
$ws.Range("J15").Value = -0.1253231084953424
$ws.Range("K15").Value = -0.3352267436446591

# Row 16: recompute I16, add J16
$ws.Range("I16").Value = 0.2284633975843539
$ws.Range("J16").Value = 0.01855976243503714

# Row 17: recompute H17, add I17
$ws.Range("H17").Value = 0.08028600715190851
$ws.Range("I17").Value = -0.1296176279974082

# Row 18: recompute G18, add H18
$ws.Range("G18").Value = -0.07715998185224648
$ws.Range("H18").Value = -0.2870636170015632

# Row 19: recompute F19, add G19
$ws.Range("F19").Value = 0.4234994746738243
$ws.Range("G19").Value = 0.2135958395245076

# Row 20: recompute E20, add F20
$ws.Range("E20").Value = 0.1431415941383551
$ws.Range("F20").Value = -0.06676204101096155

# Row 21: recompute D21, add E21
$ws.Range("D21").Value = 0.3151164519833668
$ws.Range("E21").Value = 0.1052128168340501

# Row 22: recompute C22, add D22
$ws.Range("C22").Value = 0.009253912237035311
$ws.Range("D22").Value = -0.2006497229122814

# Row 23: recompute B23, add C23
$ws.Range("B23").Value = 0.6215838649243215
$ws.Range("C23").Value = 0.4116802297750048

# Row 24: add B24
$ws.Range("B24").Value = -0.2766911554241067
